$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "Walfaanaa Magarsaa" (row 2); remaining rows shift up.
$ws.Rows.Item(2).Delete()

# Match the author's final selection state after the edit.
[void]$ws.Range("B12").Select()
